# The document contains five paragraphs of the form:
#   <id>          (Courier New, color 7f6000, run 1)
#   p103v_N       (default font,  color 000000, run 2)
#   </id>         (Courier New, color 7f6000, run 3)
# split across three separate runs. The edit merges each of these
# triples into a single run "<id>p103v_N</id>" that keeps the
# Courier-New / 7f6000 formatting of the surrounding tag runs.
#
# Doing a plain Find/Replace of the already-merged text onto itself
# causes Word to re-flow the three runs into a single run, adopting the
# formatting of the first run in the found range - exactly matching the
# target XML.

$d = $word.ActiveDocument

for ($i = 1; $i -le 5; $i++) {
    $tag = "<id>p103v_$i</id>"
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($tag, $false, $false, $false, $false, $false, $true, 1, $false, $tag, 2) | Out-Null
}
